$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.695.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.991.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.48%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.987.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.482.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.88%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.985.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "58.533.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "425.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.690"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.106"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.947"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0687"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "385.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0353"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.671.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.244"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.10%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.110"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.23%  "
